$d = $word.ActiveDocument

$d.Content.Find.Execute("4+39=43", $true, $false, $false, $false, $false, $true, 1, $false, "90-86=4", 2) | Out-Null
$d.Content.Find.Execute("93-38=55", $true, $false, $false, $false, $false, $true, 1, $false, "21+40=61", 2) | Out-Null
$d.Content.Find.Execute("40+41=81", $true, $false, $false, $false, $false, $true, 1, $false, "11+3=14", 2) | Out-Null
$d.Content.Find.Execute("86+7=93", $true, $false, $false, $false, $false, $true, 1, $false, "91-19=72", 2) | Out-Null
$d.Content.Find.Execute("26-12=14", $true, $false, $false, $false, $false, $true, 1, $false, "7+57=64", 2) | Out-Null
$d.Content.Find.Execute("39+44=83", $true, $false, $false, $false, $false, $true, 1, $false, "22-14=8", 2) | Out-Null
$d.Content.Find.Execute("18+80=98", $true, $false, $false, $false, $false, $true, 1, $false, "53+35=88", 2) | Out-Null
$d.Content.Find.Execute("52+22=74", $true, $false, $false, $false, $false, $true, 1, $false, "58+19=77", 2) | Out-Null
$d.Content.Find.Execute("75+22=97", $true, $false, $false, $false, $false, $true, 1, $false, "42-27=15", 2) | Out-Null
$d.Content.Find.Execute("84-69=15", $true, $false, $false, $false, $false, $true, 1, $false, "25+65=90", 2) | Out-Null
$d.Content.Find.Execute("63+11=74", $true, $false, $false, $false, $false, $true, 1, $false, "26+13=39", 2) | Out-Null
$d.Content.Find.Execute("28-14=14", $true, $false, $false, $false, $false, $true, 1, $false, "45-11=34", 2) | Out-Null
$d.Content.Find.Execute("8+70=78", $true, $false, $false, $false, $false, $true, 1, $false, "99-87=12", 2) | Out-Null
$d.Content.Find.Execute("8+85=93", $true, $false, $false, $false, $false, $true, 1, $false, "12+75=87", 2) | Out-Null
$d.Content.Find.Execute("66-9=57", $true, $false, $false, $false, $false, $true, 1, $false, "51+14=65", 2) | Out-Null
$d.Content.Find.Execute("53-10=43", $true, $false, $false, $false, $false, $true, 1, $false, "97-10=87", 2) | Out-Null
$d.Content.Find.Execute("32+26=58", $true, $false, $false, $false, $false, $true, 1, $false, "22+50=72", 2) | Out-Null
$d.Content.Find.Execute("33+47=80", $true, $false, $false, $false, $false, $true, 1, $false, "38+19=57", 2) | Out-Null
$d.Content.Find.Execute("28+29=57", $true, $false, $false, $false, $false, $true, 1, $false, "95-8=87", 2) | Out-Null
$d.Content.Find.Execute("20+63=83", $true, $false, $false, $false, $false, $true, 1, $false, "46+35=81", 2) | Out-Null
$d.Content.Find.Execute("19+23=42", $true, $false, $false, $false, $false, $true, 1, $false, "55+32=87", 2) | Out-Null
$d.Content.Find.Execute("25-8=17", $true, $false, $false, $false, $false, $true, 1, $false, "52+5=57", 2) | Out-Null
$d.Content.Find.Execute("22+62=84", $true, $false, $false, $false, $false, $true, 1, $false, "58-8=50", 2) | Out-Null
$d.Content.Find.Execute("99-85=14", $true, $false, $false, $false, $false, $true, 1, $false, "24+31=55", 2) | Out-Null
$d.Content.Find.Execute("63-14=49", $true, $false, $false, $false, $false, $true, 1, $false, "8+86=94", 2) | Out-Null
$d.Content.Find.Execute("47-27=20", $true, $false, $false, $false, $false, $true, 1, $false, "18+61=79", 2) | Out-Null
$d.Content.Find.Execute("98-38=60", $true, $false, $false, $false, $false, $true, 1, $false, "26+9=35", 2) | Out-Null
$d.Content.Find.Execute("6+29=35", $true, $false, $false, $false, $false, $true, 1, $false, "36-31=5", 2) | Out-Null
$d.Content.Find.Execute("71+23=94", $true, $false, $false, $false, $false, $true, 1, $false, "99-61=38", 2) | Out-Null
$d.Content.Find.Execute("69+10=79", $true, $false, $false, $false, $false, $true, 1, $false, "98-55=43", 2) | Out-Null
$d.Content.Find.Execute("61-6=55", $true, $false, $false, $false, $false, $true, 1, $false, "91+3=94", 2) | Out-Null
$d.Content.Find.Execute("63-38=25", $true, $false, $false, $false, $false, $true, 1, $false, "75-64=11", 2) | Out-Null
$d.Content.Find.Execute("48+43=91", $true, $false, $false, $false, $false, $true, 1, $false, "58-10=48", 2) | Out-Null
$d.Content.Find.Execute("11+49=60", $true, $false, $false, $false, $false, $true, 1, $false, "57+20=77", 2) | Out-Null
$d.Content.Find.Execute("99-16=83", $true, $false, $false, $false, $false, $true, 1, $false, "33+63=96", 2) | Out-Null
$d.Content.Find.Execute("18+10=28", $true, $false, $false, $false, $false, $true, 1, $false, "77-52=25", 2) | Out-Null
$d.Content.Find.Execute("49-48=1", $true, $false, $false, $false, $false, $true, 1, $false, "45+13=58", 2) | Out-Null
$d.Content.Find.Execute("60-22=38", $true, $false, $false, $false, $false, $true, 1, $false, "76+15=91", 2) | Out-Null
$d.Content.Find.Execute("93-32=61", $true, $false, $false, $false, $false, $true, 1, $false, "31+42=73", 2) | Out-Null
$d.Content.Find.Execute("80-6=74", $true, $false, $false, $false, $false, $true, 1, $false, "1+93=94", 2) | Out-Null
$d.Content.Find.Execute("49-20=29", $true, $false, $false, $false, $false, $true, 1, $false, "7-3=4", 2) | Out-Null
$d.Content.Find.Execute("2+34=36", $true, $false, $false, $false, $false, $true, 1, $false, "75-69=6", 2) | Out-Null
$d.Content.Find.Execute("72-2=70", $true, $false, $false, $false, $false, $true, 1, $false, "14+0=14", 2) | Out-Null
$d.Content.Find.Execute("70-68=2", $true, $false, $false, $false, $false, $true, 1, $false, "20+9=29", 2) | Out-Null
$d.Content.Find.Execute("29-26=3", $true, $false, $false, $false, $false, $true, 1, $false, "12+76=88", 2) | Out-Null
$d.Content.Find.Execute("9+55=64", $true, $false, $false, $false, $false, $true, 1, $false, "21-11=10", 2) | Out-Null
$d.Content.Find.Execute("62-38=24", $true, $false, $false, $false, $false, $true, 1, $false, "45+47=92", 2) | Out-Null
$d.Content.Find.Execute("64-28=36", $true, $false, $false, $false, $false, $true, 1, $false, "72-0=72", 2) | Out-Null
$d.Content.Find.Execute("32-26=6", $true, $false, $false, $false, $false, $true, 1, $false, "9+25=34", 2) | Out-Null
$d.Content.Find.Execute("17-5=12", $true, $false, $false, $false, $false, $true, 1, $false, "39+14=53", 2) | Out-Null
$d.Content.Find.Execute("63-1=62", $true, $false, $false, $false, $false, $true, 1, $false, "61-8=53", 2) | Out-Null
$d.Content.Find.Execute("76-68=8", $true, $false, $false, $false, $false, $true, 1, $false, "5+24=29", 2) | Out-Null
$d.Content.Find.Execute("26+1=27", $true, $false, $false, $false, $false, $true, 1, $false, "22-0=22", 2) | Out-Null
$d.Content.Find.Execute("13-9=4", $true, $false, $false, $false, $false, $true, 1, $false, "39+26=65", 2) | Out-Null
$d.Content.Find.Execute("91-79=12", $true, $false, $false, $false, $false, $true, 1, $false, "31-15=16", 2) | Out-Null
$d.Content.Find.Execute("54-26=28", $true, $false, $false, $false, $false, $true, 1, $false, "10-9=1", 2) | Out-Null
$d.Content.Find.Execute("16+68=84", $true, $false, $false, $false, $false, $true, 1, $false, "52-31=21", 2) | Out-Null
$d.Content.Find.Execute("55+20=75", $true, $false, $false, $false, $false, $true, 1, $false, "46+12=58", 2) | Out-Null
$d.Content.Find.Execute("74-51=23", $true, $false, $false, $false, $false, $true, 1, $false, "96-58=38", 2) | Out-Null
$d.Content.Find.Execute("46+20=66", $true, $false, $false, $false, $false, $true, 1, $false, "59+31=90", 2) | Out-Null
$d.Content.Find.Execute("86-25=61", $true, $false, $false, $false, $false, $true, 1, $false, "54-50=4", 2) | Out-Null
$d.Content.Find.Execute("92-36=56", $true, $false, $false, $false, $false, $true, 1, $false, "77-39=38", 2) | Out-Null
$d.Content.Find.Execute("99-22=77", $true, $false, $false, $false, $false, $true, 1, $false, "77-16=61", 2) | Out-Null
$d.Content.Find.Execute("69-1=68", $true, $false, $false, $false, $false, $true, 1, $false, "78-20=58", 2) | Out-Null
$d.Content.Find.Execute("63+26=89", $true, $false, $false, $false, $false, $true, 1, $false, "83-5=78", 2) | Out-Null
$d.Content.Find.Execute("22+64=86", $true, $false, $false, $false, $false, $true, 1, $false, "79-45=34", 2) | Out-Null
$d.Content.Find.Execute("79+6=85", $true, $false, $false, $false, $false, $true, 1, $false, "70+22=92", 2) | Out-Null
$d.Content.Find.Execute("94-17=77", $true, $false, $false, $false, $false, $true, 1, $false, "55+44=99", 2) | Out-Null
$d.Content.Find.Execute("97-51=46", $true, $false, $false, $false, $false, $true, 1, $false, "86-51=35", 2) | Out-Null
$d.Content.Find.Execute("44+31=75", $true, $false, $false, $false, $false, $true, 1, $false, "97+2=99", 2) | Out-Null
$d.Content.Find.Execute("3+73=76", $true, $false, $false, $false, $false, $true, 1, $false, "48+6=54", 2) | Out-Null
$d.Content.Find.Execute("85-34=51", $true, $false, $false, $false, $false, $true, 1, $false, "68-32=36", 2) | Out-Null
$d.Content.Find.Execute("58-4=54", $true, $false, $false, $false, $false, $true, 1, $false, "26-15=11", 2) | Out-Null
$d.Content.Find.Execute("43-1=42", $true, $false, $false, $false, $false, $true, 1, $false, "35+39=74", 2) | Out-Null
$d.Content.Find.Execute("46+17=63", $true, $false, $false, $false, $false, $true, 1, $false, "28+15=43", 2) | Out-Null
$d.Content.Find.Execute("40+44=84", $true, $false, $false, $false, $false, $true, 1, $false, "89-28=61", 2) | Out-Null
$d.Content.Find.Execute("94-59=35", $true, $false, $false, $false, $false, $true, 1, $false, "70-25=45", 2) | Out-Null
$d.Content.Find.Execute("57-26=31", $true, $false, $false, $false, $false, $true, 1, $false, "29+43=72", 2) | Out-Null
$d.Content.Find.Execute("25+23=48", $true, $false, $false, $false, $false, $true, 1, $false, "39+22=61", 2) | Out-Null
$d.Content.Find.Execute("28-18=10", $true, $false, $false, $false, $false, $true, 1, $false, "62-10=52", 2) | Out-Null
$d.Content.Find.Execute("50+36=86", $true, $false, $false, $false, $false, $true, 1, $false, "18+4=22", 2) | Out-Null
$d.Content.Find.Execute("66-6=60", $true, $false, $false, $false, $false, $true, 1, $false, "27-8=19", 2) | Out-Null
$d.Content.Find.Execute("8+45=53", $true, $false, $false, $false, $false, $true, 1, $false, "69+12=81", 2) | Out-Null
$d.Content.Find.Execute("72-45=27", $true, $false, $false, $false, $false, $true, 1, $false, "18-14=4", 2) | Out-Null
$d.Content.Find.Execute("96-90=6", $true, $false, $false, $false, $false, $true, 1, $false, "86-25=61", 2) | Out-Null
$d.Content.Find.Execute("57-10=47", $true, $false, $false, $false, $false, $true, 1, $false, "75-38=37", 2) | Out-Null
$d.Content.Find.Execute("98-92=6", $true, $false, $false, $false, $false, $true, 1, $false, "83-39=44", 2) | Out-Null
$d.Content.Find.Execute("77-22=55", $true, $false, $false, $false, $false, $true, 1, $false, "95-40=55", 2) | Out-Null
$d.Content.Find.Execute("61-54=7", $true, $false, $false, $false, $false, $true, 1, $false, "25-2=23", 2) | Out-Null
$d.Content.Find.Execute("21+2=23", $true, $false, $false, $false, $false, $true, 1, $false, "81-28=53", 2) | Out-Null
$d.Content.Find.Execute("11+73=84", $true, $false, $false, $false, $false, $true, 1, $false, "15+66=81", 2) | Out-Null
$d.Content.Find.Execute("39+6=45", $true, $false, $false, $false, $false, $true, 1, $false, "59-23=36", 2) | Out-Null
$d.Content.Find.Execute("95-58=37", $true, $false, $false, $false, $false, $true, 1, $false, "57+12=69", 2) | Out-Null
$d.Content.Find.Execute("80-63=17", $true, $false, $false, $false, $false, $true, 1, $false, "68+16=84", 2) | Out-Null
$d.Content.Find.Execute("21-13=8", $true, $false, $false, $false, $false, $true, 1, $false, "33+23=56", 2) | Out-Null
$d.Content.Find.Execute("51-37=14", $true, $false, $false, $false, $false, $true, 1, $false, "12+52=64", 2) | Out-Null
$d.Content.Find.Execute("50-46=4", $true, $false, $false, $false, $false, $true, 1, $false, "67-23=44", 2) | Out-Null
$d.Content.Find.Execute("52+15=67", $true, $false, $false, $false, $false, $true, 1, $false, "37-30=7", 2) | Out-Null
$d.Content.Find.Execute("38+38=76", $true, $false, $false, $false, $false, $true, 1, $false, "74-43=31", 2) | Out-Null
$d.Content.Find.Execute("74-41=33", $true, $false, $false, $false, $false, $true, 1, $false, "24+6=30", 2) | Out-Null
